$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-28 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-29 Monday", 2)
$d.Content.Find.Execute("36÷6=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "74÷9=8, 2", 2)
$d.Content.Find.Execute("40÷3=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=4, 4", 2)
$d.Content.Find.Execute("15÷9=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "93÷9=10, 3", 2)
$d.Content.Find.Execute("68÷2=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "43÷4=10, 3", 2)
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "95÷9=10, 5", 2)
$d.Content.Find.Execute("46÷4=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "33÷6=5, 3", 2)
$d.Content.Find.Execute("94÷9=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "33÷4=8, 1", 2)
$d.Content.Find.Execute("47÷4=11, 3", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=35, 1", 2)
$d.Content.Find.Execute("37÷2=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "27÷9=3, 0", 2)
$d.Content.Find.Execute("39÷9=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2)
$d.Content.Find.Execute("56÷7=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "44÷2=22, 0", 2)
$d.Content.Find.Execute("28÷4=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=26, 0", 2)
$d.Content.Find.Execute("72÷9=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=19, 2", 2)
$d.Content.Find.Execute("50÷3=16, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷2=48, 0", 2)
$d.Content.Find.Execute("79÷5=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=5, 2", 2)
$d.Content.Find.Execute("19÷9=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "99÷4=24, 3", 2)
$d.Content.Find.Execute("29÷6=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "86÷9=9, 5", 2)
$d.Content.Find.Execute("83÷4=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=21, 0", 2)
$d.Content.Find.Execute("72÷6=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2)
$d.Content.Find.Execute("75÷6=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "37÷5=7, 2", 2)
$d.Content.Find.Execute("68÷5=13, 3", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=9, 7", 2)
$d.Content.Find.Execute("70÷4=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "37÷3=12, 1", 2)
$d.Content.Find.Execute("99÷2=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=5, 4", 2)
$d.Content.Find.Execute("16÷3=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "13÷2=6, 1", 2)
$d.Content.Find.Execute("35÷9=3, 8", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=3, 2", 2)
